$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 5 new days of COVID-19 USA stats (rows 30-34) ---
# date (col A, same "d-mmm" style as existing rows), total_cases (col B)
$ws.Range("A30").Value = 43919
$ws.Range("B30").Value = 143491

$ws.Range("A31").Value = 43920
$ws.Range("B31").Value = 163788

$ws.Range("A32").Value = 43921
$ws.Range("B32").Value = 188530

$ws.Range("A33").Value = 43922
$ws.Range("B33").Value = 215003

$ws.Range("A34").Value = 43923
$ws.Range("B34").Value = 244877

# Match date formatting used by the rest of column A (style index 1, "d-mmm")
$ws.Range("A30:A34").NumberFormat = "d-mmm"

# Drift the new_cases / growth_factor formulas down through row 34, mirroring
# the existing shared formulas in columns C (B-B) and D (C/C)
$ws.Range("C30:C34").Formula = "=B30-B29"
$ws.Range("D30:D34").Formula = "=C30/C29"

# Update the view: scrolled down with E31 as the active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("E31").Select()
